$d = $word.ActiveDocument

# Locate the paragraph whose entire text is the literal placeholder "<Condition>"
# (Range.Text includes the trailing paragraph mark, hence TrimEnd()).
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd() -eq "<Condition>") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    # Range covering just the run's text, excluding the trailing paragraph mark,
    # so later the bookmark wraps only the run (not the paragraph mark too).
    $runRange = $d.Range($target.Range.Start, $target.Range.End - 1)

    # Replace the placeholder text with the new placeholder text.
    $runRange.Text = "<Assessment of Significance>"

    # Wrap the replaced run with a bookmark, as in the target revision.
    # $runRange automatically grows to cover the newly-inserted text.
    $d.Bookmarks.Add("_Hlk118981492", $runRange)
}
